$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 2.26
$ws.Cells.Item(2, 7).Value = 2.28
$ws.Cells.Item(2, 8).Value = 3.75
$ws.Cells.Item(2, 9).Value = 3.8
$ws.Cells.Item(2, 10).Value = 3.35
$ws.Cells.Item(2, 11).Value = 3.4
$ws.Cells.Item(2, 12).Value = 2.22
$ws.Cells.Item(2, 13).Value = 1.11
$ws.Cells.Item(2, 14).Value = 2.74
$ws.Cells.Item(2, 15).Value = 1.55
$ws.Cells.Item(2, 16).Value = 1.55
$ws.Cells.Item(2, 17).Value = 2.76
$ws.Cells.Item(2, 18).Value = 1.2
$ws.Cells.Item(2, 19).Value = 5.5
$ws.Cells.Item(2, 20).Value = 2.2
$ws.Cells.Item(2, 21).Value = 1.74
$ws.Cells.Item(2, 22).Value = 1.36
$ws.Cells.Item(2, 23).Value = 1.78
$ws.Cells.Item(2, 24).Value = 10.5
$ws.Cells.Item(2, 25).Value = 10.5
$ws.Cells.Item(2, 26).Value = 25
$ws.Cells.Item(2, 27).Value = 85
$ws.Cells.Item(2, 28).Value = 7.2
$ws.Cells.Item(2, 29).Value = 7.4
$ws.Cells.Item(2, 31).Value = 75
$ws.Cells.Item(2, 32).Value = 11.5
$ws.Cells.Item(2, 33).Value = 12
$ws.Cells.Item(2, 34).Value = 27
$ws.Cells.Item(2, 35).Value = 100
$ws.Cells.Item(2, 36).Value = 28
$ws.Cells.Item(2, 37).Value = 34
$ws.Cells.Item(2, 38).Value = 1000
$ws.Cells.Item(2, 39).Value = 250
$ws.Cells.Item(2, 40).Value = 29
$ws.Cells.Item(2, 41).Value = 95

# Row 3
$ws.Cells.Item(3, 6).Value = 1.77
$ws.Cells.Item(3, 7).Value = 1.84
$ws.Cells.Item(3, 8).Value = 4.5
$ws.Cells.Item(3, 9).Value = 5.1
$ws.Cells.Item(3, 11).Value = 4.6
$ws.Cells.Item(3, 14).Value = 5.1
$ws.Cells.Item(3, 16).Value = 2.38
$ws.Cells.Item(3, 17).Value = 1.65
$ws.Cells.Item(3, 18).Value = 1.55
$ws.Cells.Item(3, 19).Value = 2.6
$ws.Cells.Item(3, 20).Value = 1.68
$ws.Cells.Item(3, 21).Value = 2.44
$ws.Cells.Item(3, 22).Value = 1.25
$ws.Cells.Item(3, 23).Value = 2.18
$ws.Cells.Item(3, 28).Value = 12
$ws.Cells.Item(3, 29).Value = 10.5
$ws.Cells.Item(3, 31).Value = 120
$ws.Cells.Item(3, 32).Value = 18
$ws.Cells.Item(3, 35).Value = 250
$ws.Cells.Item(3, 36).Value = 38
$ws.Cells.Item(3, 40).Value = 9.199999999999999

# Row 4
$ws.Cells.Item(4, 6).Value = 2.3
$ws.Cells.Item(4, 7).Value = 2.4
$ws.Cells.Item(4, 8).Value = 3.6
$ws.Cells.Item(4, 9).Value = 4
$ws.Cells.Item(4, 10).Value = 3.1
$ws.Cells.Item(4, 14).Value = 3.15
$ws.Cells.Item(4, 16).Value = 1.7
$ws.Cells.Item(4, 20).Value = 1.91
$ws.Cells.Item(4, 21).Value = 1.92
$ws.Cells.Item(4, 22).Value = 1.34
$ws.Cells.Item(4, 23).Value = 1.71
$ws.Cells.Item(4, 26).Value = 27
$ws.Cells.Item(4, 27).Value = 190
$ws.Cells.Item(4, 28).Value = 8.800000000000001
$ws.Cells.Item(4, 30).Value = 16.5
$ws.Cells.Item(4, 31).Value = 55
$ws.Cells.Item(4, 32).Value = 14.5
$ws.Cells.Item(4, 33).Value = 12
$ws.Cells.Item(4, 34).Value = 21
$ws.Cells.Item(4, 35).Value = 120
$ws.Cells.Item(4, 36).Value = 34
$ws.Cells.Item(4, 37).Value = 30
$ws.Cells.Item(4, 40).Value = 26
$ws.Cells.Item(4, 41).Value = 65

# Row 5
$ws.Cells.Item(5, 6).Value = 1.51
$ws.Cells.Item(5, 8).Value = 7.8
$ws.Cells.Item(5, 9).Value = 9.6
$ws.Cells.Item(5, 11).Value = 4.8
$ws.Cells.Item(5, 14).Value = 3.85
$ws.Cells.Item(5, 16).Value = 1.96
$ws.Cells.Item(5, 17).Value = 1.96
$ws.Cells.Item(5, 18).Value = 1.35
$ws.Cells.Item(5, 19).Value = 3.5
$ws.Cells.Item(5, 21).Value = 1.81
$ws.Cells.Item(5, 22).Value = 1.12
$ws.Cells.Item(5, 28).Value = 9
$ws.Cells.Item(5, 29).Value = 19
$ws.Cells.Item(5, 40).Value = 27

# Row 6
$ws.Cells.Item(6, 7).Value = 14
$ws.Cells.Item(6, 8).Value = 1.26
$ws.Cells.Item(6, 9).Value = 1.31
$ws.Cells.Item(6, 10).Value = 6.2
$ws.Cells.Item(6, 11).Value = 7.4

# Row 7
$ws.Cells.Item(7, 8).Value = 3.25
$ws.Cells.Item(7, 9).Value = 3.65
$ws.Cells.Item(7, 10).Value = 3.05
$ws.Cells.Item(7, 11).Value = 3.5
$ws.Cells.Item(7, 12).Value = 1.46
$ws.Cells.Item(7, 14).Value = 3.35
$ws.Cells.Item(7, 15).Value = 1.37
$ws.Cells.Item(7, 17).Value = 2.12
$ws.Cells.Item(7, 18).Value = 1.3
$ws.Cells.Item(7, 19).Value = 3.9
$ws.Cells.Item(7, 20).Value = 1.82
$ws.Cells.Item(7, 21).Value = 1.98
$ws.Cells.Item(7, 25).Value = 23
$ws.Cells.Item(7, 28).Value = 9.800000000000001
$ws.Cells.Item(7, 30).Value = 30
$ws.Cells.Item(7, 32).Value = 34
$ws.Cells.Item(7, 33).Value = 23
$ws.Cells.Item(7, 36).Value = 900

# Row 8
$ws.Cells.Item(8, 12).Value = 1.47
$ws.Cells.Item(8, 14).Value = 3.65
$ws.Cells.Item(8, 15).Value = 1.36
$ws.Cells.Item(8, 16).Value = 1.88
$ws.Cells.Item(8, 17).Value = 2.1
$ws.Cells.Item(8, 18).Value = 1.34
$ws.Cells.Item(8, 19).Value = 3.9
$ws.Cells.Item(8, 20).Value = 1.88
$ws.Cells.Item(8, 21).Value = 2.1
$ws.Cells.Item(8, 23).Value = 1.79
$ws.Cells.Item(8, 25).Value = 13
$ws.Cells.Item(8, 28).Value = 9
$ws.Cells.Item(8, 30).Value = 15.5
$ws.Cells.Item(8, 31).Value = 42
$ws.Cells.Item(8, 32).Value = 13
$ws.Cells.Item(8, 34).Value = 18
$ws.Cells.Item(8, 38).Value = 40
$ws.Cells.Item(8, 39).Value = 110
$ws.Cells.Item(8, 40).Value = 19
$ws.Cells.Item(8, 41).Value = 50

# Row 9
$ws.Cells.Item(9, 8).Value = 6.4
$ws.Cells.Item(9, 9).Value = 6.6
$ws.Cells.Item(9, 10).Value = 4.7
$ws.Cells.Item(9, 11).Value = 4.8
$ws.Cells.Item(9, 16).Value = 2.36
$ws.Cells.Item(9, 17).Value = 1.72
$ws.Cells.Item(9, 18).Value = 1.52
$ws.Cells.Item(9, 21).Value = 2.2
$ws.Cells.Item(9, 29).Value = 10
$ws.Cells.Item(9, 30).Value = 23
$ws.Cells.Item(9, 31).Value = 75
$ws.Cells.Item(9, 41).Value = 95

# Row 10
$ws.Cells.Item(10, 6).Value = 1.92
$ws.Cells.Item(10, 7).Value = 2.02
$ws.Cells.Item(10, 8).Value = 3.65
$ws.Cells.Item(10, 9).Value = 4.1
$ws.Cells.Item(10, 10).Value = 4
$ws.Cells.Item(10, 14).Value = 5.4
$ws.Cells.Item(10, 16).Value = 2.52
$ws.Cells.Item(10, 18).Value = 1.61
$ws.Cells.Item(10, 21).Value = 2.42
$ws.Cells.Item(10, 23).Value = 1.98
$ws.Cells.Item(10, 30).Value = 17.5
$ws.Cells.Item(10, 34).Value = 16.5
$ws.Cells.Item(10, 40).Value = 9.4

# Row 11
$ws.Cells.Item(11, 6).Value = 2.34
$ws.Cells.Item(11, 7).Value = 2.38
$ws.Cells.Item(11, 12).Value = 1.42
$ws.Cells.Item(11, 14).Value = 4
$ws.Cells.Item(11, 21).Value = 2.24
$ws.Cells.Item(11, 23).Value = 1.73
$ws.Cells.Item(11, 35).Value = 44

# Row 12
$ws.Cells.Item(12, 6).Value = 8.4
$ws.Cells.Item(12, 7).Value = 8.6
$ws.Cells.Item(12, 8).Value = 1.44
$ws.Cells.Item(12, 9).Value = 1.45
$ws.Cells.Item(12, 12).Value = 1.36
$ws.Cells.Item(12, 13).Value = 1.05
$ws.Cells.Item(12, 14).Value = 4.9
$ws.Cells.Item(12, 15).Value = 1.24
$ws.Cells.Item(12, 16).Value = 2.32
$ws.Cells.Item(12, 17).Value = 1.73
$ws.Cells.Item(12, 18).Value = 1.52
$ws.Cells.Item(12, 19).Value = 2.86
$ws.Cells.Item(12, 20).Value = 1.97
$ws.Cells.Item(12, 22).Value = 3.2
$ws.Cells.Item(12, 23).Value = 1.13
$ws.Cells.Item(12, 25).Value = 8.800000000000001
$ws.Cells.Item(12, 26).Value = 8.4
$ws.Cells.Item(12, 27).Value = 12
$ws.Cells.Item(12, 36).Value = 270
$ws.Cells.Item(12, 39).Value = 130
$ws.Cells.Item(12, 40).Value = 190
$ws.Cells.Item(12, 41).Value = 6.2
